$wb = $excel.ActiveWorkbook

# Reset the view of the currently-active sheet (FukuokaKousokuDouroYomikata)
# before adding + activating the new sheet, matching the diff (tabSelected=false,
# selection moves to B1).
$wsFukuoka = $wb.Worksheets.Item("FukuokaKousokuDouroYomikata")
$wsFukuoka.Activate()
$wsFukuoka.Range("B1").Select()

# Add the new sheet right after FukuokaKousokuDouroYomikata (i.e. as the last sheet)
$newSheet = $wb.Worksheets.Add($null, $wsFukuoka)
$newSheet.Name = "ZenkokuKousokuDouroYomikata"

# Header row
$newSheet.Cells.Item(1, 1).Value = "kanji"
$newSheet.Cells.Item(1, 2).Value = "hiragana"

# Data rows: kanji place name / hiragana reading for the all-Japan (Hokkaido)
# expressway reading dictionary
$data = @(
    @("国縫", "くんぬい"),
    @("長万部", "おしゃまんべ"),
    @("静狩", "しずかり"),
    @("豊浦", "とようら"),
    @("豊浦噴火湾", "とようらふんかわん"),
    @("虻田洞爺湖", "あぶたとうやこ"),
    @("伊達", "だて"),
    @("有珠山", "うすざん"),
    @("室蘭", "むろらん"),
    @("本輪西", "もとわにし"),
    @("登別室蘭", "のぼりべつむろらん"),
    @("富浦", "とみうら"),
    @("登別東", "のぼりべつひがし"),
    @("萩野", "はぎの"),
    @("白老", "しらおい"),
    @("樽前", "たるまえ"),
    @("苫小牧西", "とまこまいにし"),
    @("苫小牧東", "とまこまいひがし"),
    @("美沢", "みさわ"),
    @("千歳", "ちとせ"),
    @("千歳恵庭", "ちとせえにわ"),
    @("恵庭", "えにわ"),
    @("輪厚", "わっつ"),
    @("北広島", "きたひろしま"),
    @("札幌南", "さっぽろみなみ"),
    @("大谷地", "おおやち"),
    @("北郷", "きたごう"),
    @("札幌", "さっぽろ"),
    @("江別西", "えべつにし"),
    @("野幌", "のっぽろ"),
    @("江別東", "えべつひがし"),
    @("岩見沢", "いわみざわ"),
    @("三笠", "みかさ"),
    @("美唄", "びばい"),
    @("茶志内", "ちゃしない"),
    @("奈井江砂川", "ないえすながわ"),
    @("砂川", "すながわ"),
    @("滝川", "たきかわ"),
    @("深川", "ふかがわ"),
    @("音江", "おとえ"),
    @("旭川鷹栖", "あさひかわたかす"),
    @("旭川北", "あさひかわきた"),
    @("比布大雪", "ひっぷだいせつ"),
    @("和寒", "わさむ"),
    @("札樽道", "さっそんどう"),
    @("小樽", "おたる"),
    @("朝里", "あさり"),
    @("銭函", "ぜにばこ"),
    @("金山", "かなやま"),
    @("手稲", "ていね"),
    @("札幌西", "さっぽろにし"),
    @("新川", "しんかわ"),
    @("札幌北", "さっぽろきた"),
    @("伏古", "ふしこ"),
    @("雁来", "かりき"),
    @("千歳東", "ちとせひがし"),
    @("追分町", "おいわけちょう"),
    @("夕張", "ゆうばり"),
    @("十勝清水", "とかちしみず"),
    @("十勝平原", "とかちへいげん"),
    @("芽室", "めむろ"),
    @("音更帯広", "おとふけおびひろ"),
    @("池田", "いけだ"),
    @("深川西", "ふかがわにし")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $newSheet.Cells.Item($row, 1).Value = $data[$i][0]
    $newSheet.Cells.Item($row, 2).Value = $data[$i][1]
}

# Column widths (approximate the target character widths)
$newSheet.Columns.Item(1).ColumnWidth = 10.61
$newSheet.Columns.Item(2).ColumnWidth = 17.85

# Make the new sheet the active tab + set cursor/selection, matching the diff
# (tabSelected=true, activeCell=A2, activeTab index = 14)
$newSheet.Activate()
$newSheet.Range("A2").Select()
